# Generate Report for Handback
# Adds a new handback row (file 4cad4159-f8ff-4c93-9907-e530f54578c4.md) to
# the Overview, zh-cn and de-de sheets, mirroring the existing rows for
# 0caed798-38b1-4c7b-85e3-ba430893ab93.md / 468b9ed0-1b47-45f9-a922-0e8b8714b07e.md.

$wb = $excel.ActiveWorkbook

$fileGuid   = "4cad4159-f8ff-4c93-9907-e530f54578c4"
$mdName     = "$fileGuid.md"
$mdPath     = "e2e\$fileGuid.md"
$statusSync = "Handed back: in sync with en-US"
$ext        = ".md"

$zhXlf  = "$fileGuid.dcf83435d3d741171ade37365af8bd5f636a3b30.zh-cn.xlf"
$deXlf  = "$fileGuid.dcf83435d3d741171ade37365af8bd5f636a3b30.de-de.xlf"

$mdDate      = "2016-08-29 08:47:52"
$zhHoDate    = "2016-08-29 08:47:47"
$zhHbDate    = "2016-08-29 08:48:20"
$deHoDate    = "2016-08-29 08:47:52"
$deHbDate    = "2016-08-29 08:48:27"

$mdUrlBase     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c6e8a1f4b9c2d7e5a0f18b6c4d9e2a7f5b8c1d3/e2e/$mdName"
$zhCnMdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7a1d4f9c2b8e6a3d5f0c17b9e4a6d2c8f1b5a3e7/e2e/$mdName"
$deDeMdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2f8b5d1a9c4e7b0d3a6f19c5e8b2d4a7f0c3e6b9/e2e/$mdName"

# ---------------------------------------------------------------------
# Sheet "Overview" -> new row 4
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4, 1).Value = $mdName
$wsOverview.Cells.Item(4, 2).Value = $mdPath
$wsOverview.Cells.Item(4, 3).Value = $ext
$wsOverview.Cells.Item(4, 5).Value = $statusSync
$wsOverview.Cells.Item(4, 6).Value = $statusSync
$wsOverview.Cells.Item(4, 7).Value = $mdDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $mdUrlBase, "", "", $mdPath)

$overviewTable = $wsOverview.ListObjects.Item(1)
$overviewTable.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> new row 4
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(4, 1).Value  = $mdName
$wsZhCn.Cells.Item(4, 2).Value  = $ext
$wsZhCn.Cells.Item(4, 3).Value  = $statusSync
$wsZhCn.Cells.Item(4, 4).Value  = "e2e"
$wsZhCn.Cells.Item(4, 5).Value  = "ht"
$wsZhCn.Cells.Item(4, 6).Value  = "True"
$wsZhCn.Cells.Item(4, 7).Value  = $zhXlf
$wsZhCn.Cells.Item(4, 8).Value  = $zhHoDate
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(4, 9).Value  = $mdName
$wsZhCn.Cells.Item(4, 10).Value = $zhXlf
$wsZhCn.Cells.Item(4, 11).Value = $zhHbDate
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(4, 12).Value = ""
$wsZhCn.Cells.Item(4, 13).Value = "True"
$wsZhCn.Cells.Item(4, 14).Value = ""
$wsZhCn.Cells.Item(4, 15).Value = "False"
$wsZhCn.Cells.Item(4, 16).Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $mdUrlBase, "", "", $mdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), $zhCnMdUrl, "", "", $mdName)

$zhCnTable = $wsZhCn.ListObjects.Item(1)
$zhCnTable.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet "de-de" -> new row 4
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(4, 1).Value  = $mdName
$wsDeDe.Cells.Item(4, 2).Value  = $ext
$wsDeDe.Cells.Item(4, 3).Value  = $statusSync
$wsDeDe.Cells.Item(4, 4).Value  = "e2e"
$wsDeDe.Cells.Item(4, 5).Value  = "ht"
$wsDeDe.Cells.Item(4, 6).Value  = "True"
$wsDeDe.Cells.Item(4, 7).Value  = $deXlf
$wsDeDe.Cells.Item(4, 8).Value  = $deHoDate
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(4, 9).Value  = $mdName
$wsDeDe.Cells.Item(4, 10).Value = $deXlf
$wsDeDe.Cells.Item(4, 11).Value = $deHbDate
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(4, 12).Value = ""
$wsDeDe.Cells.Item(4, 13).Value = "True"
$wsDeDe.Cells.Item(4, 14).Value = ""
$wsDeDe.Cells.Item(4, 15).Value = "False"
$wsDeDe.Cells.Item(4, 16).Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $mdUrlBase, "", "", $mdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), $deDeMdUrl, "", "", $mdName)

$deDeTable = $wsDeDe.ListObjects.Item(1)
$deDeTable.Resize($wsDeDe.Range("A1:P4"))
